{"js": "// Remove the leftover yellow highlight from two ToDo bullets:\n//  1) \"Would having the staff lines extend beyond the advance width ...\"\n//  2) \"A b-flat on position A ... (see Liber Usualis 1681).\"\n// Every run in those paragraphs (and the paragraph mark itself) carried\n// <w:highlight w:val=\"yellow\"/>; clearing Paragraph.font.highlightColor\n// removes it from the whole paragraph in one go.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needles = [\n  \"Would having the staff lines extend beyond the advance width\",\n  \"A b-flat on position A\"\n];\n\nconst targets = paragraphs.items.filter((p) =>\n  needles.some((needle) => p.text.indexOf(needle) !== -1)\n);\n\nfor (const paragraph of targets) {\n  paragraph.font.highlightColor = null;\n}\n\nawait context.sync();\n", "ps1": "# Remove the leftover yellow highlight from two ToDo bullets:\n#  1) \"Would having the staff lines extend beyond the advance width ...\"\n#  2) \"A b-flat on position A ... (see Liber Usualis 1681).\"\n# Every run of text in those paragraphs was marked\n# Range.HighlightColorIndex = wdYellow (7); clear it back to\n# wdNoHighlight (0).\n\n$d = $word.ActiveDocument\n\n$needles = @(\n    \"Would having the staff lines extend beyond the advance width\",\n    \"A b-flat on position A\"\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text\n\n    $isTarget = $false\n    foreach ($needle in $needles) {\n        if ($text -like \"*$needle*\") {\n            $isTarget = $true\n        }\n    }\n\n    if ($isTarget) {\n        $para.Range.HighlightColorIndex = 0\n    }\n}\n"}
